# Fruta / hortaliza, semanal
# Insert two new weekly price rows (2022-07-11, serial 44753) at the top of the
# data block for "Terminal Hortofrutícola Agro Chillán - Limón" (row 682),
# pushing the existing historical rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows before the current row 682, shifting rows 682:713 down to 684:715
$ws.Rows("682:683").Insert()

# Row 682: Limón "1a amarillo"
$ws.Range("A682").Value2 = 7
$ws.Range("B682").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C682").Value2 = "Ñuble"
$ws.Range("D682").Value2 = 44753
$ws.Range("D682").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E682").Value2 = 16
$ws.Range("F682").Value2 = "Fruta"
$ws.Range("G682").Value2 = 100102
$ws.Range("H682").Value2 = "Cítricos"
$ws.Range("I682").Value2 = 100102003
$ws.Range("J682").Value2 = "Limón"
$ws.Range("K682").Value2 = "Sin especificar"
$ws.Range("L682").Value2 = "1a amarillo"
$ws.Range("M682").Value2 = 120
$ws.Range("N682").Value2 = 4000
$ws.Range("O682").Value2 = 4500
$ws.Range("P682").Value2 = 4250
$ws.Range("Q682").Value2 = "$/malla 16 kilos"
$ws.Range("R682").Value2 = "Región de O'Higgins"
$ws.Range("S682").Value2 = 266
$ws.Range("T682").Value2 = 16

# Row 683: Limón "2a amarillo"
$ws.Range("A683").Value2 = 7
$ws.Range("B683").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C683").Value2 = "Ñuble"
$ws.Range("D683").Value2 = 44753
$ws.Range("D683").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E683").Value2 = 16
$ws.Range("F683").Value2 = "Fruta"
$ws.Range("G683").Value2 = 100102
$ws.Range("H683").Value2 = "Cítricos"
$ws.Range("I683").Value2 = 100102003
$ws.Range("J683").Value2 = "Limón"
$ws.Range("K683").Value2 = "Sin especificar"
$ws.Range("L683").Value2 = "2a amarillo"
$ws.Range("M683").Value2 = 60
$ws.Range("N683").Value2 = 3500
$ws.Range("O683").Value2 = 3500
$ws.Range("P683").Value2 = 3500
$ws.Range("Q683").Value2 = "$/malla 16 kilos"
$ws.Range("R683").Value2 = "Región de O'Higgins"
$ws.Range("S683").Value2 = 219
$ws.Range("T683").Value2 = 16
